$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row -> hashtable of column letter -> new value.
# Only the cells that actually changed (per the diff) are listed.
$updates = [ordered]@{
    2  = @{ D = '29.428.70'; E = '  -0.49%  ' }
    3  = @{ D = '1.851.24';  E = '  -0.20%  ' }
    4  = @{ E = '  -0.04%  ' }
    5  = @{ D = '241.26';    E = '  -1.09%  ' }
    6  = @{ E = '  -1.56%  ' }
    7  = @{ E = '  +0.01%  ' }
    8  = @{ D = '4.788.15';  E = '  +152.84%  ' }
    9  = @{ D = '4.918.54';  E = '  +126.05%  ' }
    10 = @{ D = '0.07576';   E = '  +0.60%  ' }
    11 = @{ D = '0.2961';    E = '  -1.64%  ' }
    12 = @{ D = '24.61';     E = '  +0.98%  ' }
    13 = @{ D = '0.07728';   E = '  +0.75%  ' }
    14 = @{ D = '5.000';     E = '  -1.14%  ' }
    15 = @{ D = '0.6846';    E = '  -0.86%  ' }
    16 = @{ D = '83.04';     E = '  -1.22%  ' }
    17 = @{ D = '0.000009929'; E = '  +3.43%  ' }
    18 = @{ E = '  -1.18%  ' }
    19 = @{ D = '29.458.62'; E = '  -0.56%  ' }
    20 = @{ D = '231.91';    E = '  -2.75%  ' }
    21 = @{ E = '  -1.24%  ' }
    22 = @{ E = '  +0.02%  ' }
    23 = @{ D = '7.614';     E = '  -1.40%  ' }
    24 = @{ E = '  -0.04%  ' }
    25 = @{ D = '156.04';    E = '  -0.80%  ' }
    26 = @{ D = '0.1389';    E = '  -1.99%  ' }
    27 = @{ D = '8.409';     E = '  -1.56%  ' }
    28 = @{ D = '17.71';     E = '  -0.77%  ' }
    29 = @{ D = '4.970.10';  E = '  +139.74%  ' }
    30 = @{ D = '1.471';     E = '  -1.20%  ' }
    31 = @{ D = '0.05757';   E = '  -3.69%  ' }
    32 = @{ D = '1.259';     E = '  -0.12%  ' }
    33 = @{ D = '4.133';     E = '  -0.47%  ' }
    34 = @{ D = '4.022';     E = '  -1.41%  ' }
    35 = @{ D = '1.857';     E = '  -1.68%  ' }
    36 = @{ E = '  -1.36%  ' }
    37 = @{ D = '0.7166';    E = '  -1.08%  ' }
    38 = @{ E = '  -0.29%  ' }
    39 = @{ B = 'FraxShare';         C = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs';            D = '7.128';    E = '  +15.13%  ' }
    40 = @{ B = 'Maker';             C = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr';            D = '1.255.12'; E = '  +3.48%  ' }
    41 = @{ B = 'MXToken';           C = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx';            D = '2.802';    E = '  +0.65%  ' }
    42 = @{ B = 'VeChain';           C = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet';          D = '0.01806';  E = '  +1.41%  ' }
    43 = @{ B = 'TrustWalletToken';  C = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt';     D = '0.9052';   E = '  -1.30%  ' }
    44 = @{ E = '  -0.02%  ' }
    45 = @{ D = '101.77';    E = '  -0.29%  ' }
    46 = @{ D = '67.00';     E = '  -0.81%  ' }
    47 = @{ D = '7.155';     E = '  -3.12%  ' }
    48 = @{ D = '9.191';     E = '  -0.03%  ' }
    49 = @{ E = '  -1.10%  ' }
    50 = @{ D = '1.686';     E = '  +1.28%  ' }
    51 = @{ B = 'Frax';              C = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax';                D = '1.096';    E = '  +9.64%  ' }
}

foreach ($row in $updates.Keys) {
    $cells = $updates[$row]
    if ($cells.ContainsKey('B')) {
        $ws.Cells.Item($row, 2).Value = $cells['B']
    }
    if ($cells.ContainsKey('C')) {
        $ws.Cells.Item($row, 3).Value = $cells['C']
    }
    if ($cells.ContainsKey('D')) {
        # Column D stores price strings as text (e.g. "29.428.70"); force text
        # format so Excel does not reinterpret the value as a number and
        # strip formatting (trailing zeros, multiple '.' separators, etc.).
        $cell = $ws.Cells.Item($row, 4)
        $cell.NumberFormat = "@"
        $cell.Value = $cells['D']
    }
    if ($cells.ContainsKey('E')) {
        $ws.Cells.Item($row, 5).Value = $cells['E']
    }
}
